$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # The underlying data series for "Cutoff" (col B) / "Reaction_number" (col C)
    # was recomputed starting 4 steps later than before. Pull the values that
    # used to live 4 rows further down and write them into rows 2-16, leaving
    # column A's 0..14 index untouched, then drop the now-stale trailing rows.
    for ($r = 2; $r -le 16; $r++) {
        $srcRow = $r + 4
        $bVal = $ws.Cells.Item($srcRow, 2).Value2
        $cVal = $ws.Cells.Item($srcRow, 3).Value2
        $ws.Cells.Item($r, 2).Value = $bVal
        $ws.Cells.Item($r, 3).Value = $cVal
    }

    $ws.Range("A17:C20").EntireRow.Delete()
}
